$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# Rebuild the border styling of the merged-header row (B1:D1): the middle
# cell (C1) gets a top+bottom border, the last cell (D1) gets a
# top+bottom+right border (these match border definitions already present
# in the workbook's style table: borderId 4 and 5 respectively).
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders(8).LineStyle = 1
$ws1.Range("C1").Borders(9).LineStyle = 1

$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders(8).LineStyle = 1
$ws1.Range("D1").Borders(10).LineStyle = 1
$ws1.Range("D1").Borders(9).LineStyle = 1

# Anonymize the "fedcore" label
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# Re-use the exact styles just built on sheet 1 (instead of re-deriving
# the borders step by step again), so every one of the four header cells
# across both sheets (C1/D1 on sheet1, C1/D1/F1/G1 on sheet2) maps onto
# the same two new style entries with no left-over intermediate styles.
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

# Anonymize the "fedcore" labels
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
